$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to 2-decimal "custom accuracy" rounded figures
$ws.Range("B5").Value = 8.17
$ws.Range("C5").Value = 5.79
$ws.Range("D5").Value = 0.83
$ws.Range("E5").Value = 17.72
$ws.Range("F5").Value = 14.17
$ws.Range("G5").Value = 6.43
$ws.Range("H5").Value = 26.96
$ws.Range("I5").Value = 9.89
$ws.Range("J5").Value = 4.29
$ws.Range("K5").Value = 6.22
$ws.Range("L5").Value = 7.11
$ws.Range("M5").Value = 7.44
$ws.Range("N5").Value = 2.06
$ws.Range("O5").Value = 6.39
$ws.Range("P5").Value = 9.01
$ws.Range("Q5").Value = 5.57
$ws.Range("R5").Value = 0.74
$ws.Range("S5").Value = 0.48
$ws.Range("T5").Value = 89.96
$ws.Range("U5").Value = 17.96
$ws.Range("V5").Value = 5.9
$ws.Range("W5").Value = 11.88
$ws.Range("X5").Value = 6.18
$ws.Range("Y5").Value = 1.16
$ws.Range("Z5").Value = 12.79
$ws.Range("AA5").Value = 5.21
$ws.Range("AB5").Value = 4.73
$ws.Range("AC5").Value = 5.54
$ws.Range("AD5").Value = 7.42
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 24.65
$ws.Range("AG5").Value = 3.23
$ws.Range("AH5").Value = 7.38

# Remove row 6 entirely (data was trimmed to fewer rows)
$ws.Rows("6").Delete()
